# Auto-generated edit script applying F-column ('想去人数' / interested-count) updates
# as described by the commit 'Update gh-pages to output generated at 456a3b4'.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 33   # F2: 28 -> 33
$ws.Cells.Item(3, 6).Value = 51   # F3: 50 -> 51
$ws.Cells.Item(4, 6).Value = 951   # F4: 948 -> 951
$ws.Cells.Item(5, 6).Value = 1239   # F5: 1237 -> 1239
$ws.Cells.Item(6, 6).Value = 1690   # F6: 1685 -> 1690
$ws.Cells.Item(7, 6).Value = 900   # F7: 901 -> 900
$ws.Cells.Item(8, 6).Value = 561   # F8: 559 -> 561
$ws.Cells.Item(9, 6).Value = 2443   # F9: 2425 -> 2443
$ws.Cells.Item(10, 6).Value = 705   # F10: 699 -> 705
$ws.Cells.Item(11, 6).Value = 566   # F11: 562 -> 566
$ws.Cells.Item(13, 6).Value = 8   # F13: 6 -> 8
$ws.Cells.Item(15, 6).Value = 325   # F15: 323 -> 325
$ws.Cells.Item(16, 6).Value = 201   # F16: 198 -> 201
$ws.Cells.Item(18, 6).Value = 2100   # F18: 2098 -> 2100
$ws.Cells.Item(19, 6).Value = 1226   # F19: 1224 -> 1226
$ws.Cells.Item(20, 6).Value = 695   # F20: 693 -> 695
$ws.Cells.Item(22, 6).Value = 2602   # F22: 2601 -> 2602
$ws.Cells.Item(24, 6).Value = 22   # F24: 21 -> 22
$ws.Cells.Item(25, 6).Value = 512   # F25: 510 -> 512
$ws.Cells.Item(27, 6).Value = 321   # F27: 317 -> 321
$ws.Cells.Item(28, 6).Value = 1756   # F28: 1750 -> 1756
$ws.Cells.Item(30, 6).Value = 527   # F30: 525 -> 527
$ws.Cells.Item(34, 6).Value = 4537   # F34: 4533 -> 4537
$ws.Cells.Item(35, 6).Value = 105   # F35: 94 -> 105

# Sheet 2: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 405   # F2: 404 -> 405
$ws.Cells.Item(3, 6).Value = 382   # F3: 381 -> 382
$ws.Cells.Item(7, 6).Value = 41   # F7: 40 -> 41
$ws.Cells.Item(14, 6).Value = 306   # F14: 304 -> 306
$ws.Cells.Item(26, 6).Value = 237   # F26: 228 -> 237
$ws.Cells.Item(27, 6).Value = 9   # F27: 8 -> 9
$ws.Cells.Item(28, 6).Value = 248   # F28: 247 -> 248
$ws.Cells.Item(34, 6).Value = 3   # F34: 2 -> 3
$ws.Cells.Item(37, 6).Value = 63   # F37: 62 -> 63

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 1408   # F4: 1407 -> 1408
$ws.Cells.Item(6, 6).Value = 504   # F6: 500 -> 504
$ws.Cells.Item(7, 6).Value = 171   # F7: 165 -> 171

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 1408   # F3: 1407 -> 1408
$ws.Cells.Item(5, 6).Value = 504   # F5: 500 -> 504
$ws.Cells.Item(6, 6).Value = 382   # F6: 381 -> 382
$ws.Cells.Item(7, 6).Value = 33   # F7: 28 -> 33
$ws.Cells.Item(8, 6).Value = 51   # F8: 50 -> 51
$ws.Cells.Item(9, 6).Value = 951   # F9: 948 -> 951
$ws.Cells.Item(10, 6).Value = 1239   # F10: 1237 -> 1239
$ws.Cells.Item(11, 6).Value = 1690   # F11: 1685 -> 1690
$ws.Cells.Item(12, 6).Value = 41   # F12: 40 -> 41
$ws.Cells.Item(15, 6).Value = 900   # F15: 901 -> 900
$ws.Cells.Item(16, 6).Value = 561   # F16: 559 -> 561
$ws.Cells.Item(17, 6).Value = 2443   # F17: 2425 -> 2443
$ws.Cells.Item(18, 6).Value = 705   # F18: 699 -> 705
$ws.Cells.Item(19, 6).Value = 566   # F19: 562 -> 566
$ws.Cells.Item(22, 6).Value = 325   # F22: 323 -> 325
$ws.Cells.Item(24, 6).Value = 201   # F24: 198 -> 201
$ws.Cells.Item(27, 6).Value = 2100   # F27: 2098 -> 2100
$ws.Cells.Item(28, 6).Value = 1226   # F28: 1224 -> 1226
$ws.Cells.Item(29, 6).Value = 695   # F29: 693 -> 695
$ws.Cells.Item(32, 6).Value = 2602   # F32: 2601 -> 2602
$ws.Cells.Item(35, 6).Value = 22   # F35: 21 -> 22
$ws.Cells.Item(36, 6).Value = 512   # F36: 510 -> 512
$ws.Cells.Item(38, 6).Value = 171   # F38: 165 -> 171
$ws.Cells.Item(41, 6).Value = 1756   # F41: 1750 -> 1756
$ws.Cells.Item(42, 6).Value = 237   # F42: 228 -> 237
$ws.Cells.Item(43, 6).Value = 527   # F43: 525 -> 527
$ws.Cells.Item(47, 6).Value = 4537   # F47: 4533 -> 4537
$ws.Cells.Item(48, 6).Value = 105   # F48: 94 -> 105
$ws.Cells.Item(49, 6).Value = 63   # F49: 62 -> 63
